# PRACTICA 5: Cambios en los ficheros y memoria casi acabada
#
# Inserts two new model rows ("Perceptron con capa oculta (500) ..." and
# "Perceptron con capa oculta (750) ...") right after the existing
# "(100) 'relu'" row (new rows 13 and 14), pushing the remainder of the
# results table down by two rows. Also tidies up the sheet view (scrolled
# back to the top, with C15 selected) and narrows column A, which no
# longer needs to "best fit" the longest label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows for the two new models, right before the old
#     row 13 ("Perceptron con 2 capa oculta (10, 10) ..."). Restrict the
#     insert to columns A:C so the "SOBREAJUSTE" annotation living in D53
#     is not dragged down along with everything else; it will be
#     re-anchored to D53 explicitly afterwards.
$ws.Range("A13:C14").Insert()

$ws.Range("A13").Value = "Perceptron con capa oculta (500) 'relu' y salida 'softmax'"
$ws.Range("B13").Value = 93.86
$ws.Range("C13").Value = 93.87

$ws.Range("A14").Value = "Perceptron con capa oculta (750) 'relu' y salida 'softmax'"
$ws.Range("B14").Value = 94.12
$ws.Range("C14").Value = 94.13

# --- The D53 "SOBREAJUSTE" cell got pushed to D55 by the row shift above;
#     move it back up so it keeps sitting next to the same data row it
#     originally annotated (which is now two rows further down, at row 55).
$ws.Range("D55").Cut($ws.Range("D53"))
$ws.Range("D55").Clear()

# --- Column A no longer needs to auto-fit the longest label.
$ws.Columns.Item(1).ColumnWidth = 47.6

# --- Reset the view: scroll back to the top-left corner and select C15.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C15").Select()
